{"js": "// The second paragraph currently reads (split across runs, with a\n// \"_GoBack\" bookmark sitting in the middle of the bold file-extension\n// list):\n//   ... files (DOCX, DOC, PDF, HTML, XPS, R[_GoBack]TF and TXT) from ...\n// It should read, as a single contiguous bold run with the bookmark\n// removed:\n//   ... files (DOCX, DOC, PDF, HTML, XPS, RTF and TXT) from ...\n\nconst body = context.document.body;\n\n// search() can match text that spans multiple runs/bookmarks, so look\n// for the full phrase and replace it in one shot - this merges the\n// \"R\" run and the \"TF and TXT\" run back into a single run.\nconst results = body.search(\"XPS, RTF and TXT\", { matchCase: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"XPS, RTF and TXT\", \"Replace\");\n}\n\n// Remove the now-stray \"_GoBack\" bookmark left over in the middle of\n// the (now merged) run.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The second paragraph currently reads (split across runs, with a\n# \"_GoBack\" bookmark sitting in the middle of the bold file-extension\n# list):\n#   ... files (DOCX, DOC, PDF, HTML, XPS, R[_GoBack]TF and TXT) from ...\n# It should read, as a single contiguous bold run with the bookmark\n# removed:\n#   ... files (DOCX, DOC, PDF, HTML, XPS, RTF and TXT) from ...\n\n$d = $word.ActiveDocument\n\n# Find & Replace across the whole bold phrase. Word's Find engine can\n# match text that spans multiple runs/bookmarks, and replacing it\n# collapses everything it matched (including the stray \"_GoBack\"\n# bookmark sitting between the two runs) into freshly written run(s).\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"XPS, RTF and TXT\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"XPS, RTF and TXT\"\n$find.Execute([ref]$null, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$null, 2) | Out-Null\n\n# Safety net in case the bookmark survived the replace for some reason\n# (it shouldn't - Find/Replace over a range removes bookmarks fully\n# contained in it - but make sure it is gone either way).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n"}
